# "Generate Report for Archive"
#
# The localization status for this item moved on from handoff, so the
# "Ready for handoff" status text becomes "In Translation" everywhere it
# is used (Overview!E2:F2, and the Status column - C2 - on each of the
# per-locale report sheets). The status column on each sheet is then
# narrowed to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status cells -----------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column -----------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column -----------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
